$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (test user) cleanup ---
# forename/surname/street_nr/city cleared out
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("G2").Value = ""
# zip reset to 0
$ws.Range("F2").Value = 0
# email set to placeholder
$ws.Range("H2").Value = "..._...@...."

# --- Row 3 (max_mustermann) isLoggedIn flipped off ---
$ws.Range("K3").Value = $false

# --- Selection moved to M3 ---
$ws.Range("M3").Select() | Out-Null
